$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A318").Value = "IMX-USD"
$ws.Range("A319").Value = "TAO-USD"
$ws.Range("A320").Value = "MNT-USD"
$ws.Range("A321").Value = "GRT-USD"
